# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
# The "Date" column (BF) held a mangled string like "2-24-2011-12";
# correct it to the proper ISO-ish date string "2012-02-24" for every
# data row, keeping the cells as plain text (not auto-converted dates).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$correctDate = "2012-02-24"

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Range("BF$row")
    # Force text interpretation so Excel doesn't reinterpret the
    # "yyyy-mm-dd"-looking string as a real date value/serial number.
    $cell.NumberFormat = "@"
    $cell.Value = $correctDate
    # Restore the default "Normal" style so no stray per-cell style
    # (beyond the text format) lingers on these cells, matching the
    # original (unstyled) cells.
    $cell.Style = "Normal"
}
